# "Generate Report for Handback"
#
# A new handoff/handback round-trip completed for the
# "81d9b0b3-ef76-46a3-98eb-92e1d16c015c.md" file (the 2nd data row, row 3,
# in both the zh-cn and de-de report sheets). Refresh the report with the
# newly recorded timestamps:
#   - zh-cn  : Correspond Handoff Datetime (H3) / Correspond Handback DateTime (K3)
#   - de-de  : Correspond Handoff Datetime (H3) / Correspond Handback DateTime (K3)
# and roll the de-de "Latest HO Xliff Generate Date" shown on the Overview
# sheet (G3) forward to match the new de-de handoff timestamp.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# zh-cn: row 3 is the 81d9b0b3-...md file
$wsZhCn.Range("H3").Value = "2016-08-31 02:55:17"
$wsZhCn.Range("K3").Value = "2016-08-31 02:55:33"

# de-de: row 3 is the 81d9b0b3-...md file
$wsDeDe.Range("H3").Value = "2016-08-31 02:55:21"
$wsDeDe.Range("K3").Value = "2016-08-31 02:55:41"

# Overview: de-de "Latest HO Xliff Generate Date" for the same file (row 3)
$wsOverview.Range("G3").Value = "2016-08-31 02:55:21"
